$d = $word.ActiveDocument
$p8 = $d.Paragraphs.Item(8)
$startPoint = $p8.Range.Start
$newRange = $d.Range($startPoint, $startPoint)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Corpo"/><w:spacing w:before="1" w:after="160" w:line="259" w:lineRule="auto"/></w:pPr></w:p>'
$newRange.InsertXML($xml)
